# Apply the data/content changes described by the commit
# ("new example with real inp file included") to the mapping sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the numeric mapping values (columns C, D, E) ---
# Row 2 (S1)
$ws.Range("C2").Value = 0
$ws.Range("E2").Value = 1

# Row 3 (S2)
$ws.Range("C3").Value = -1
$ws.Range("D3").Value = 0.8

# Row 4 (S3)
$ws.Range("C4").Value = -1
$ws.Range("D4").Value = 0.8

# Row 5 (S4)
$ws.Range("C5").Value = -2
$ws.Range("E5").Value = 1

# --- Column widths: two new custom-width columns (D, E) appear, and the
#     best-fit widths of C, G, H shift slightly to match the new content. ---
$ws.Columns.Item(3).ColumnWidth = 13.1666666666667
$ws.Columns.Item(4).ColumnWidth = 9.5
$ws.Columns.Item(5).ColumnWidth = 9.66666666666667
$ws.Columns.Item(7).ColumnWidth = 10.6666666666667
$ws.Columns.Item(8).ColumnWidth = 13.6666666666667

# --- Selection moves from G5 to E5 ---
$ws.Range("E5").Select()
